# Add "2022-Q4" data: a new sheet inserted right after "总计" (position 1),
# plus a new row at the top of the "总计" summary table.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert a brand-new worksheet right after "总计" and name it "2022-Q4".
#    NOTE: worksheet references resolved via Item(index)/Item(name) are
#    live/positional in this engine - any Add() shifts every sheet after
#    the insertion point, so we must re-resolve sheets *after* the Add()
#    rather than keep using a reference captured beforehand.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)
$newSheet = $wb.Worksheets.Add($null, $total)
$newSheet.Name = "2022-Q4"

# Re-resolve "2022-Q3" AFTER the insert - it has been pushed to index 3.
$q3 = $wb.Worksheets.Item(3)

# Copy header row formatting (bold/border header style) and the A-column
# index style from "2022-Q3" so the new sheet matches the look of the
# other quarterly sheets. Column A on row 1 is intentionally left blank
# (it is blank on every quarterly sheet), so copy B1:H1 and A2:H2 only.
$q3.Range("B1:H1").Copy($newSheet.Range("B1:H1"))
$q3.Range("A2:H2").Copy($newSheet.Range("A2:H2"))

# Overwrite row 2 with the actual 2022-Q4 fund data. Numeric-looking
# values (fund code / size / position figures) are entered quote-prefixed
# so they stay text (matching the source data), matching column C's
# (plain-text) cell style look by stripping the quote-prefix marker style
# that Value-assignment adds automatically.
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "'002810"
$newSheet.Range("C2").Value = "金信转型创新成长灵活配置混合"
$newSheet.Range("D2").Value = "'3.84"
$newSheet.Range("E2").Value = "'89.18"
$newSheet.Range("F2").Value = "'7.32"
$newSheet.Range("G2").Value = "'0.2811"
$newSheet.Range("H2").Value = 2

$plainStyle = $newSheet.Range("C2").Style
$newSheet.Range("B2").Style = $plainStyle
$newSheet.Range("D2").Style = $plainStyle
$newSheet.Range("E2").Style = $plainStyle
$newSheet.Range("F2").Style = $plainStyle
$newSheet.Range("G2").Style = $plainStyle

# ---------------------------------------------------------------------
# 2. Update the "总计" sheet: push the existing quarter rows (2..7) down
#    by one row (6..7 shift to 7..8, ..., 2 shifts to 3), then write the
#    new "2022-Q4" entry into row 2. Re-sequence the A-column index.
# ---------------------------------------------------------------------
for ($r = 7; $r -ge 2; $r--) {
    $srcRange = "A" + $r + ":D" + $r
    $dstRow = $r + 1
    $dstRange = "A" + $dstRow + ":D" + $dstRow
    $total.Range($srcRange).Copy($total.Range($dstRange))
}

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 1
$total.Range("D2").Value = 0.28

for ($r = 3; $r -le 8; $r++) {
    $total.Range("A" + $r).Value = $r - 2
}
